# Fix Implementation Approach timeline formatting: remove markdown italic
# asterisk syntax around the phase month ranges and merge them into the
# bold phase-title run, e.g.
#   "Phase 1: Discovery & Assessment" + " *(Months 1-2)*"
#   -> "Phase 1: Discovery & Assessment (Months 1-2)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

$fixes = @{
    1 = "Phase 1: Discovery & Assessment (Months 1-2)"
    5 = "Phase 2: Migration Execution (Months 3-6)"
    9 = "Phase 3: Optimization (Months 7-9)"
}

foreach ($idx in $fixes.Keys) {
    $para = $tr.Paragraphs($idx, 1)
    $run1 = $para.Runs(1, 1)
    $run2 = $para.Runs(2, 1)
    $run1.Text = $fixes[$idx]
    $run2.Text = ""
}
